$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5318679213523865
$ws.Range("B1").Value = 1.446839451789856
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 1.533384084701538
$ws.Range("E1").Value = 1.441953897476196
